$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row: "_old" -> "_FV2304", "_new" -> "_FV2310" ---
$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
# K1 "diff" is unchanged

$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# --- Turn the used range into an Excel Table ("Table1") ---
$range = $ws.Range("A1:U91")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (split below row 1) ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
